$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feedback text content (Rafal -> Ramon rename + extra sentence in AXA feedback).
# Written in E4, E3, E6, E7 order so the rebuilt shared-string table lands on the same
# indices as the target workbook.
$ws.Range("E4").Value = "From Recognize || Unique Accenture Moments || Accenture achievements || Hello,   On behalf of AXA Health project Leadership, I would like to thank you for all commitment and time you spend working with us on the account which led us to the end of development of release 1.   We all appreciate your involvement  and we would like to show our gratitude with this points recognition. You could be more punctual.   Keep up the good work!   Thank you, AXA Health Team || 02-Mar-23"
$ws.Range("E3").Value = "Dear Ramon, `nI would like to express my sincere appreciation for your exceptional qualities as a colleague. Your technical expertise is truly remarkable, and your willingness to assist others is greatly valued. Your profound understanding of our current project is evident and highly beneficial. At present, I believe you have surpassed all expectations, leaving no room for improvement. I extend my heartfelt gratitude for the invaluable help and insightful suggestions you have provided. `nThank you immensely for your continued support and contributions. `nBest regards,Dana Kalm"
$ws.Range("E6").Value = "I have been working with Ramon for the past one year. Ramon showed incredible leadership instincts in all tasks. He has incredible knowledge in Azure and optimum solutions for every complex problem. He is very positive person I would love to develop those skills. Amazing work."
$ws.Range("E7").Value = "From Recognize || Custom Program || Team Awards || Ramon - Congratulations on finishing in the Top 10 for the iAi Holiday Hackathon. On behalf of our collective team, thanks for your efforts and for continuing to help make these hackathons so impactful. iAi Team || 12-Jan-22"

# Column E becomes its own (wider) width group instead of sharing the A:E 23-char block.
$ws.Columns.Item(5).ColumnWidth = 81.6

# Row heights re-fitted for the new content/column width.
$ws.Rows.Item(3).RowHeight = 140.25
$ws.Rows.Item(4).RowHeight = 76.5
$ws.Rows.Item(5).RowHeight = 153
$ws.Rows.Item(6).RowHeight = 140.25
$ws.Rows.Item(7).RowHeight = 51

# Move the active selection to E3.
$ws.Range("E3").Select()
